$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.691.75'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '2.097.33'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').Value = "'1.009"
$ws.Range('E4').Value = '  +0.59%  '
$ws.Range('D5').Value = "'343.26"
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').Value = "'0.5165"
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'0.4377"
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('D9').Value = "'53.55"
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').Value = "'0.09184"
$ws.Range('E10').Value = '  +2.45%  '
$ws.Range('D11').Value = "'1.166"
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').Value = "'24.61"
$ws.Range('E12').Value = '  -5.34%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'6.759"
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '2.049.16'
$ws.Range('E14').Value = '  -2.06%  '
$ws.Range('D15').Value = "'8.142"
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('D16').Value = "'102.36"
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').Value = "'0.00001151"
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = "'1.009"
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').Value = "'0.06673"
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').Value = "'6.199"
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').Value = '29.755.91'
$ws.Range('E23').Value = '  -1.46%  '
$ws.Range('D24').Value = "'12.64"
$ws.Range('E24').Value = '  -2.10%  '
$ws.Range('E25').Value = '  -2.04%  '
$ws.Range('D26').Value = '2.335.50'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('D28').Value = "'161.90"
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').Value = "'2.488"
$ws.Range('E29').Value = '  -2.84%  '
$ws.Range('D30').Value = "'133.29"
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('E31').Value = '  -4.83%  '
$ws.Range('D32').Value = "'1.670"
$ws.Range('E32').Value = '  +1.43%  '
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('D34').Value = "'6.190"
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('D35').Value = "'3.961"
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('D36').Value = "'6.315"
$ws.Range('E36').Value = '  +6.63%  '
$ws.Range('D37').Value = "'10.44"
$ws.Range('E37').Value = '  +2.14%  '
$ws.Range('D38').Value = "'0.02578"
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').Value = "'0.06697"
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('D40').Value = "'0.6992"
$ws.Range('E40').Value = '  +1.95%  '
$ws.Range('E41').Value = '  +6.52%  '
$ws.Range('D42').Value = "'12.41"
$ws.Range('E42').Value = '  -2.38%  '
$ws.Range('D43').Value = "'0.2211"
$ws.Range('E43').Value = '  -4.73%  '
$ws.Range('D44').Value = "'0.6793"
$ws.Range('E44').Value = '  +5.68%  '
$ws.Range('D45').Value = "'14.25"
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('D46').Value = "'2.315"
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').Value = "'0.00000000361"
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('D48').Value = "'3.610"
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('E49').Value = '  +3.27%  '
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').Value = "'80.95"
$ws.Range('E51').Value = '  -3.58%  '

# Reset style on cells that required a quote-prefix to stay text,
# so they don't pick up an explicit text number format / style index.
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
